# "gold and workout updates"
# - Insert a new "Golf" sheet between "Hockey" and "Working Out", populated
#   with three rounds of golf data.
# - Append a new "Time" column + a new row to the "Poutine" sheet.
# - Add a "How Long" / "How many" column pair + three new rows to the
#   "Working Out" sheet.
# - Update various view-state bits (selected cell, active sheet/tab).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 0) Add the new "Golf" sheet first (right after "Hockey", before
#    "Working Out") so every worksheet lookup by name done afterwards
#    sees the final tab order.
# ---------------------------------------------------------------------
$hockey = $wb.Worksheets.Item("Hockey")
$golf = $wb.Worksheets.Add($null, $hockey)
$golf.Name = "Golf"

# ---------------------------------------------------------------------
# 1) Poutine sheet: insert a "Time" column before "Comment", add a row.
# ---------------------------------------------------------------------
$poutine = $wb.Worksheets.Item("Poutine")

# Shift the existing "Comment" column (G) one column right to H, and
# give the freed-up G column a new "Time" header/body. (Values are
# written literally - rather than copied via a .Value read - to avoid
# depending on Range value read-back.)
$poutine.Range("H4").Value = "Comment"
$poutine.Range("H5").Value = "bad"
$poutine.Range("H6").Value = "serviceable. Better than expected"

$poutine.Range("G4").Value = "Time"
$poutine.Range("G5").Value = "11pm"
$poutine.Range("G6").Value = "1pm"

# New row 7.
$poutine.Range("B6").Copy($poutine.Range("B7"))
$poutine.Range("B7").Value = 45457
$poutine.Range("C7").Value = "Squirleys (bar)"
$poutine.Range("D7").Value = 8.9
$poutine.Range("E7").Value = 10
$poutine.Range("F7").Value = "M"
$poutine.Range("G7").Value = "11 30pm"
$poutine.Range("H7").Value = "surpringingly good!!"

$poutine.Range("F8").Select()

# ---------------------------------------------------------------------
# 2) Working Out sheet: insert "How Long"/"How many" columns, add rows.
# ---------------------------------------------------------------------
$workout = $wb.Worksheets.Item("Working Out")

# Original header: C=Date D=Plank E=Pushup F=Headstand G=Boxing H=Yoga
# New header:      C=Date D=Plank E=HowLong F=Pushup G=HowMany H=Headstand I=Boxing J=Yoga
# Shift Yoga (H->J), Boxing (G->I), Headstand (F->H), Pushup (E->F), then
# fill in the two newly freed columns (E, G).
$workout.Range("J3").Value = "Yoga"
$workout.Range("I3").Value = "Boxing"
$workout.Range("H3").Value = "Headstand"
$workout.Range("F3").Value = "Pushup"
$workout.Range("E3").Value = "How Long"
$workout.Range("G3").Value = "How many"

# Row 4 data: Headstand ("y") moves from F4 to H4; date serial updated.
$workout.Range("C4").Value = 41072
$workout.Range("H4").Value = "y"
$workout.Range("F4").Value = ""

# New rows 5-7.
$workout.Range("C4").Copy($workout.Range("C5"))
$workout.Range("C5").Value = 41073
$workout.Range("F5").Value = "y"
$workout.Range("G5").Value = 20
$workout.Range("H5").Value = "y"
$workout.Range("I5").Value = "y"
$workout.Range("J5").Value = "y"

$workout.Range("C4").Copy($workout.Range("C6"))
$workout.Range("C6").Value = 41075
$workout.Range("H6").Value = "Y"
$workout.Range("J6").Value = "Y"

$workout.Range("C4").Copy($workout.Range("C7"))
$workout.Range("C7").Value = 41076
$workout.Range("F7").Value = "Y"
$workout.Range("G7").Value = 20
$workout.Range("H7").Value = "Y"
$workout.Range("I7").Value = "Y"
$workout.Range("J7").Value = "Y"

# ---------------------------------------------------------------------
# 3) Populate the "Golf" sheet.
# ---------------------------------------------------------------------
$golf.Columns.Item(3).ColumnWidth = $hockey.Columns.Item(2).ColumnWidth

$golf.Range("C4").Value = "Date"
$golf.Range("D4").Value = "Course"
$golf.Range("E4").Value = "Score"
$golf.Range("F4").Value = "Par"
$golf.Range("G4").Value = "Rating"
$golf.Range("H4").Value = "Slope"
$golf.Range("I4").Value = "Pars"
$golf.Range("J4").Value = "Muligans"
$golf.Range("K4").Value = "Generous Gimmies"
$golf.Range("L4").Value = "Team"
$golf.Range("M4").Value = "Comment"

$poutine.Range("B6").Copy($golf.Range("C5"))
$golf.Range("C5").Value = 45444
$golf.Range("D5").Value = "Scarlet Woods"
$golf.Range("E5").Value = 40
$golf.Range("F5").Value = 29
$golf.Range("G5").Value = 61.1
$golf.Range("H5").Value = 93
$golf.Range("I5").Value = 2
$golf.Range("J5").Value = 1
$golf.Range("K5").Value = 2
$golf.Range("M5").Value = "9 holes"

$poutine.Range("B6").Copy($golf.Range("C6"))
$golf.Range("C6").Value = 45454
$golf.Range("D6").Value = "Scarlet Woods"
$golf.Range("E6").Formula = "=37+44"
$golf.Range("F6").Value = 61
$golf.Range("G6").Value = 61.1
$golf.Range("H6").Value = 93
$golf.Range("I6").Value = 3
$golf.Range("J6").Value = 2
$golf.Range("K6").Value = 3
$golf.Range("M6").Value = "was getting duck slice need to work on backswing and followthrough. Focus on wrists and thumb"

$poutine.Range("B6").Copy($golf.Range("C7"))
$golf.Range("C7").Value = 45458
$golf.Range("D7").Value = "Centennial north"
$golf.Range("E7").Value = 45
$golf.Range("I7").Value = 1
$golf.Range("J7").Value = 3
$golf.Range("K7").Value = 3
$golf.Range("M7").Value = "played poorly but had some good shots. Went to the range after, and discovered when make DELIBERTELY slow back swing, significantly longer and better. Keep trying!"

$golf.Range("F18:F19").Select()

# ---------------------------------------------------------------------
# 4) Misc view-state touch-ups to match the saved workbook state.
# ---------------------------------------------------------------------
$movies = $wb.Worksheets.Item("Movies")
$movies.Range("D12").Select()

$videogames = $wb.Worksheets.Item("Videogames")
$videogames.Range("J23").Select()

$hockey = $wb.Worksheets.Item("Hockey")
$hockey.Range("N24").Select()

$workout = $wb.Worksheets.Item("Working Out")
$workout.Range("D7").Select()
$workout.Activate()
